$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1688.8572
$ws.Range("I19").Value = 1684.4
$ws.Range("J19").Value = 1700
$ws.Range("K19").Value = 1684.4
$ws.Range("L19").Value = 1700
$ws.Range("M19").Value = -1509.4
$ws.Range("N19").Value = -2050
$ws.Range("H80").Value = 1280.875
$ws.Range("I80").Value = 962.5
$ws.Range("K80").Value = 2887.5
$ws.Range("M80").Value = -1889.5
$ws.Range("H83").Value = 1280.875
$ws.Range("I83").Value = 962.5
$ws.Range("K83").Value = 8662.5
$ws.Range("M83").Value = -3670.5
$ws.Range("H88").Value = 1641.2858
$ws.Range("J88").Value = 1747.75
$ws.Range("L88").Value = 1747.75
$ws.Range("N88").Value = -2559.75
$ws.Range("H91").Value = 1641.2858
$ws.Range("J91").Value = 1747.75
$ws.Range("L91").Value = 1747.75
$ws.Range("N91").Value = -4555.75
$ws.Range("H113").Value = 5000
$ws.Range("J113").Value = 5000
$ws.Range("L113").Value = 5000
$ws.Range("N113").Value = -11508
$ws.Range("H132").Value = 2529.6924
$ws.Range("I132").Value = 2598.8333
$ws.Range("K132").Value = 7796.499899999999
$ws.Range("M132").Value = -5266.499899999999
$ws.Range("H137").Value = 3223.125
$ws.Range("J137").Value = 4329.3335
$ws.Range("L137").Value = 12988.0005
$ws.Range("N137").Value = -18088.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2128
$ws.Range("I45").Value = 2006
$ws.Range("K45").Value = 2006
$ws.Range("M45").Value = -1629
$ws.Range("H61").Value = 2647.6924
$ws.Range("I61").Value = 2647.6924
$ws.Range("K61").Value = 2647.6924
$ws.Range("M61").Value = -2435.6924
$ws.Range("H122").Value = 2499
$ws.Range("I122").Value = 2499
$ws.Range("K122").Value = 7497
$ws.Range("M122").Value = -5047
$ws.Range("H132").Value = 1951.2858
$ws.Range("I132").Value = 783.8
$ws.Range("K132").Value = 2351.4
$ws.Range("M132").Value = 178.6000000000004
$ws.Range("H136").Value = 2647.6924
$ws.Range("I136").Value = 2647.6924
$ws.Range("K136").Value = 7943.0772
$ws.Range("M136").Value = -5393.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").Value = ""
$ws.Range("H86").Value = 1200
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 1200
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 1200
$ws.Range("M86").Value = ""
$ws.Range("N86").Value = -3446
$ws.Range("H89").Value = 1200
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 1200
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 6000
$ws.Range("M89").Value = ""
$ws.Range("N89").Value = -17232
$ws.Range("H107").Value = 1007.8333
$ws.Range("I107").Value = 938
$ws.Range("J107").Value = 1147.5
$ws.Range("K107").Value = 938
$ws.Range("L107").Value = 1147.5
$ws.Range("M107").Value = 982
$ws.Range("N107").Value = -4987.5
$ws.Range("H123").Value = 48999
$ws.Range("J123").Value = 48999
$ws.Range("L123").Value = 48999
$ws.Range("N123").Value = -58799

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5175000
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = ""
$ws.Range("H56").Value = 52000
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").Value = ""
$ws.Range("H62").Value = 4579.8
$ws.Range("J62").Value = 5999.5
$ws.Range("L62").Value = 5999.5
$ws.Range("N62").Value = -7247.5
$ws.Range("H65").Value = 4579.8
$ws.Range("J65").Value = 5999.5
$ws.Range("L65").Value = 29997.5
$ws.Range("N65").Value = -36237.5
$ws.Range("H106").Value = 24999.5
$ws.Range("J106").Value = 24999.5
$ws.Range("L106").Value = 24999.5
$ws.Range("N106").Value = -27523.5
$ws.Range("H122").Value = 872.625
$ws.Range("I122").Value = 881.5
$ws.Range("J122").Value = 846
$ws.Range("K122").Value = 2644.5
$ws.Range("L122").Value = 2538
$ws.Range("M122").Value = -194.5
$ws.Range("N122").Value = -7438
$ws.Range("H134").Value = 4417.25
$ws.Range("I134").Value = 4441
$ws.Range("K134").Value = 13323
$ws.Range("M134").Value = -10788

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I5").Value = 1000
$ws.Range("K5").Value = 3000
$ws.Range("M5").Value = -2888
$ws.Range("H33").Value = 141
$ws.Range("J33").Value = 165
$ws.Range("L33").Value = 990
$ws.Range("N33").Value = -1556
$ws.Range("H44").Value = 1336.6666
$ws.Range("I44").Value = 1386
$ws.Range("J44").Value = 1275
$ws.Range("K44").Value = 4158
$ws.Range("L44").Value = 3825
$ws.Range("M44").Value = -3760
$ws.Range("N44").Value = -4621
$ws.Range("H68").Value = 1498.1428
$ws.Range("J68").Value = 1331.6666
$ws.Range("L68").Value = 3994.9998
$ws.Range("N68").Value = -5616.9998
$ws.Range("H71").Value = 1498.1428
$ws.Range("J71").Value = 1331.6666
$ws.Range("L71").Value = 11984.9994
$ws.Range("N71").Value = -20096.9994
$ws.Range("H80").Value = 20000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 20000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 60000
$ws.Range("M80").Value = ""
$ws.Range("N80").Value = -61872
$ws.Range("H83").Value = 20000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 20000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 180000
$ws.Range("M83").Value = ""
$ws.Range("N83").Value = -189360
$ws.Range("I135").Value = 1000
$ws.Range("K135").Value = 9000
$ws.Range("M135").Value = -6465

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 27249.5
$ws.Range("J101").Value = 27249.5
$ws.Range("L101").Value = 27249.5
$ws.Range("N101").Value = -33739.5
$ws.Range("H113").Value = 5293.25
$ws.Range("I113").Value = 4870.4
$ws.Range("K113").Value = 4870.4
$ws.Range("M113").Value = -2700.4
$ws.Range("H122").Value = 12503020
$ws.Range("I122").Value = 13891354
$ws.Range("J122").Value = 8008
$ws.Range("K122").Value = 41674062
$ws.Range("L122").Value = 24024
$ws.Range("M122").Value = -41671612
$ws.Range("N122").Value = -28924

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1270
$ws.Range("I22").Value = 1900
$ws.Range("K22").Value = 1900
$ws.Range("M22").Value = -1605
$ws.Range("H27").Value = 1270
$ws.Range("I27").Value = 1900
$ws.Range("K27").Value = 1900
$ws.Range("M27").Value = -1793
$ws.Range("H122").Value = 3492.4
$ws.Range("I122").Value = 3379.889
$ws.Range("J122").Value = 4505
$ws.Range("K122").Value = 10139.667
$ws.Range("L122").Value = 13515
$ws.Range("M122").Value = -7689.667000000001
$ws.Range("N122").Value = -18415
$ws.Range("H123").Value = 77996.5
$ws.Range("J123").Value = 77996.5
$ws.Range("L123").Value = 77996.5
$ws.Range("N123").Value = -87796.5
$ws.Range("H130").Value = 66500.5
$ws.Range("J130").Value = 66500.5
$ws.Range("L130").Value = 66500.5
$ws.Range("N130").Value = -76540.5
$ws.Range("H132").Value = 8312.5
$ws.Range("I132").Value = 11000
$ws.Range("J132").Value = 7416.6665
$ws.Range("K132").Value = 33000
$ws.Range("L132").Value = 22249.9995
$ws.Range("M132").Value = -30470
$ws.Range("N132").Value = -27309.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 6152.0527
$ws.Range("I113").Value = 9184.5
$ws.Range("K113").Value = 27553.5
$ws.Range("M113").Value = -25383.5
$ws.Range("H122").Value = 1937.4286
$ws.Range("I122").Value = 1474.5
$ws.Range("K122").Value = 4423.5
$ws.Range("M122").Value = -1973.5
$ws.Range("H132").Value = 1724.6818
$ws.Range("J132").Value = 2347.8
$ws.Range("L132").Value = 7043.400000000001
$ws.Range("N132").Value = -12103.4
$ws.Range("H133").Value = 110357.5
$ws.Range("I133").Value = 100000
$ws.Range("J133").Value = 120715
$ws.Range("K133").Value = 100000
$ws.Range("L133").Value = 120715
$ws.Range("M133").Value = -94940
$ws.Range("N133").Value = -130835
